# "put the email up top"
#
# Changes applied to the single slide of JobOpenings_SingleSlide.pptx:
#  1. Move/resize the body "Content Placeholder 2" down a bit (its Y offset
#     changes) to make room above it.
#  2. Fix a typo in the body text: "...development chops.." (double period)
#     becomes "...development chops." (single period), with the trailing
#     period ending up as its own run.
#  3. Move the "WebPlatformJobs@adobe.com" textbox up to the top of the
#     slide and make its text bigger/bold.
#  4. Move the "html.adobe.com" textbox and drop its right-alignment so it
#     reads as left aligned like the other footer text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> Points helper. PowerPoint COM measures Shape geometry in points
# (1 pt = 12700 EMU) but the canonical OOXML stores English Metric Units.
# A tiny epsilon nudges the float conversion so it rounds to the exact
# target EMU value instead of truncating one unit short.
function EmuToPt($emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

# --- 1. Content Placeholder 2: shift down -------------------------------
$body = $s.Shapes.Item("Content Placeholder 2")
$body.Left = EmuToPt(228600)
$body.Top  = EmuToPt(1382838)
$body.Width  = EmuToPt(8686800)
$body.Height = EmuToPt(4547162)

# --- 2. Fix "...chops.." -> "...chops." ---------------------------------
$bodyTr = $body.TextFrame.TextRange
$fullText = $bodyTr.Text
$idx = $fullText.IndexOf("chops..")
if ($idx -ge 0) {
    # 1-based COM character index of the two trailing periods.
    $start = $idx + 1 + ("chops").Length
    $dots = $bodyTr.Characters($start, 2)
    $dots.Text = "."
}

# --- 3. WebPlatformJobs@adobe.com textbox: move up + embiggen ----------
$email = $s.Shapes.Item("TextBox 4")
$email.Left = EmuToPt(2679931)
$email.Top  = EmuToPt(933862)
$email.Width  = EmuToPt(3781798)
$email.Height = EmuToPt(400110)

$emailTr = $email.TextFrame.TextRange
$emailTr.Font.Size = 20
$emailTr.Font.Bold = $true

# --- 4. html.adobe.com textbox: reposition + left-align -----------------
$site = $s.Shapes.Item("TextBox 5")
$site.Left = EmuToPt(74702)
$site.Top  = EmuToPt(5998339)
$site.Width  = EmuToPt(2010615)
$site.Height = EmuToPt(369332)

$siteTr = $site.TextFrame.TextRange
$siteTr.ParagraphFormat.Alignment = 1
